# Applies the "testdata: include more cases" edit to Resources.xlsx
#
# Summary of changes:
#  - Switch the active/selected tab from "classes" to "Owner"
#  - Select cell A9 on "Owner" (was B9)
#  - "classes" sheet (sheet1):
#      * widen column J to fit new, longer values
#      * pad/alter several existing J-column values with extra whitespace /
#        extra text
#      * add a bunch of new rows (16-23) full of deliberately
#        whitespace-only / "invalid" values used as additional test
#        fixtures
#  - "Owner" sheet (sheet2):
#      * pad a few existing column-A values with extra whitespace
#      * add new rows (16-19) with whitespace-only / "invalid" values

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# classes sheet
# ---------------------------------------------------------------------
$classes = $wb.Worksheets.Item("classes")

# new column J is wider to fit the longer strings now stored there
$classes.Columns.Item(10).ColumnWidth = 22.25

# existing cells in column J get extra annotation / whitespace
$classes.Cells.Item(2, 10).Value = " Resource , dcterms:fantasy   "
$classes.Cells.Item(4, 10).Value = "Resource  "
$classes.Cells.Item(5, 10).Value = "   Resource   "
$classes.Cells.Item(8, 10).Value = "StillImageRepresentation, dcterms:image  "

# row 16: "invalid" "because" "super" "is" "missing"
$classes.Cells.Item(16, 1).Value = "invalid"
$classes.Cells.Item(16, 2).Value = "because"
$classes.Cells.Item(16, 2).Font.Color = 0
$classes.Cells.Item(16, 3).Value = "super"
$classes.Cells.Item(16, 3).Font.Color = 0
$classes.Cells.Item(16, 4).Value = "is"
$classes.Cells.Item(16, 5).Value = "missing"

# row 17: "invalid" "because" "name" "is" ... "missing"
$classes.Cells.Item(17, 2).Value = "invalid"
$classes.Cells.Item(17, 3).Value = "because"
$classes.Cells.Item(17, 3).Font.Color = 0
$classes.Cells.Item(17, 4).Value = "name"
$classes.Cells.Item(17, 4).Font.Color = 0
$classes.Cells.Item(17, 5).Value = "is"
$classes.Cells.Item(17, 10).Value = "missing"

# row 18: new blank-ish cell + widen existing whitespace cell
$classes.Cells.Item(18, 1).Value = "  "
$classes.Cells.Item(18, 10).Value = "   "

# row 19: two new whitespace-only cells
$classes.Cells.Item(19, 4).Value = "   "
$classes.Cells.Item(19, 5).Value = "    "

# row 20 (new)
$classes.Cells.Item(20, 2).Value = "  "
$classes.Cells.Item(20, 3).Value = "     "
$classes.Cells.Item(20, 10).Value = "   "

# row 21 (new)
$classes.Cells.Item(21, 2).Value = "             "
$classes.Cells.Item(21, 4).Value = "  "
$classes.Cells.Item(21, 8).Value = "   "

# row 22 (new)
$classes.Cells.Item(22, 5).Value = "     "

# row 23 (new)
$classes.Cells.Item(23, 8).Value = "      "

# ---------------------------------------------------------------------
# Owner sheet
# ---------------------------------------------------------------------
$owner = $wb.Worksheets.Item("Owner")

$owner.Cells.Item(2, 1).Value = "   hasAnthroponym   "
$owner.Cells.Item(3, 1).Value = "  isOwnerOf"
$owner.Cells.Item(4, 1).Value = "correspondsToGenericAnthroponym    "

# row 16 (new)
$owner.Cells.Item(16, 1).Value = "    "
$owner.Cells.Item(16, 2).Value = "   "

# row 17 (new)
$owner.Cells.Item(17, 1).Value = "invalid"

# row 18 (new)
$owner.Cells.Item(18, 2).Value = "invalid"

# row 19 (new)
$owner.Cells.Item(19, 2).Value = "     "

# ---------------------------------------------------------------------
# Active tab / selection
# ---------------------------------------------------------------------
# Make "Owner" the active sheet (classes loses tabSelected, Owner gains
# it + workbookView gets activeTab="1"), then move its selection to A9.
$owner.Activate()
$owner.Range("A9").Select()
